$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.783.90"
$ws.Range("E2").Value = "  -2.32%  "
$ws.Range("D3").Value = "3.260.61"
$ws.Range("E3").Value = "  -1.72%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "572.15"
$ws.Range("E5").Value = "  -1.05%  "
$ws.Range("D6").Value = "176.90"
$ws.Range("E6").Value = "  -5.27%  "
$ws.Range("D7").Value = "0.624"
$ws.Range("E7").Value = "  +2.62%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.124"
$ws.Range("E9").Value = "  -4.08%  "
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("D11").Value = "0.398"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").Value = "3.834.54"
$ws.Range("E12").Value = "  -1.51%  "
$ws.Range("E13").Value = "  -3.85%  "
$ws.Range("D14").Value = "65.878.02"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "26.27"
$ws.Range("E15").Value = "  -4.36%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.267.82"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "0.0000161"
$ws.Range("E17").Value = "  -3.88%  "
$ws.Range("D18").Value = "432.79"
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "5.53"
$ws.Range("E19").Value = "  -3.14%  "
$ws.Range("D20").Value = "13.07"
$ws.Range("E20").Value = "  -3.93%  "
$ws.Range("D21").Value = "7.36"
$ws.Range("E21").Value = "  -4.89%  "
$ws.Range("D22").Value = "71.79"
$ws.Range("E22").Value = "  -3.09%  "
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").Value = "3.417.53"
$ws.Range("E24").Value = "  -1.09%  "
$ws.Range("D25").Value = "0.503"
$ws.Range("E25").Value = "  -2.94%  "
$ws.Range("E26").Value = "  +3.61%  "
$ws.Range("E27").Value = "  -6.40%  "
$ws.Range("D28").Value = "8.80"
$ws.Range("E28").Value = "  -3.11%  "
$ws.Range("D29").Value = "0.998"
$ws.Range("E29").Value = "  -0.30%  "
$ws.Range("D30").Value = "1.92"
$ws.Range("E30").Value = "  -2.66%  "
$ws.Range("D31").Value = "22.18"
$ws.Range("E31").Value = "  -3.36%  "
$ws.Range("E32").Value = "  +0.10%  "
$ws.Range("D33").Value = "5.11"
$ws.Range("E33").Value = "  -4.35%  "
$ws.Range("D34").Value = "6.55"
$ws.Range("E34").Value = "  -3.71%  "
$ws.Range("D35").Value = "1.17"
$ws.Range("E35").Value = "  -5.77%  "
$ws.Range("D36").Value = "158.12"
$ws.Range("E36").Value = "  -3.00%  "
$ws.Range("D37").Value = "1.42"
$ws.Range("E37").Value = "  -6.54%  "
$ws.Range("D38").Value = "26.37"
$ws.Range("E38").Value = "  -3.00%  "
$ws.Range("E39").Value = "  -4.21%  "
$ws.Range("D40").Value = "2.753.43"
$ws.Range("E40").Value = "  -0.61%  "
$ws.Range("D41").Value = "0.774"
$ws.Range("E41").Value = "  -2.21%  "
$ws.Range("D42").Value = "4.28"
$ws.Range("E42").Value = "  -4.60%  "
$ws.Range("D43").Value = "40.22"
$ws.Range("E43").Value = "  +0.29%  "
$ws.Range("D44").Value = "6.04"
$ws.Range("E44").Value = "  -3.29%  "
$ws.Range("D45").Value = "0.0652"
$ws.Range("E45").Value = "  -3.15%  "
$ws.Range("D46").Value = "320.62"
$ws.Range("E46").Value = "  -1.69%  "
$ws.Range("E47").Value = "  -6.14%  "
$ws.Range("D48").Value = "23.11"
$ws.Range("E48").Value = "  -6.98%  "
$ws.Range("D49").Value = "0.0265"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("E50").Value = "  +1.82%  "
$ws.Range("E51").Value = "  +0.07%  "
